$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column index, new text value.
# Columns D (Price) and E (Volume) hold numeric-looking text (e.g. "1.001",
# "29.916.35", "  +0.14%  ") that must stay plain text, matching the source
# data's inlineStr cells. Forcing NumberFormat to "@" (Text) before the
# assignment stops Excel from silently parsing them as numbers/dates, and
# resetting the style to "Normal" afterwards keeps the cell's look unchanged
# (no left-over custom number format).
$updates = @(
    @{Row=2; Col=4; Value="29.916.35"; ForceText=$True},
    @{Row=2; Col=5; Value="  +0.14%  "; ForceText=$True},
    @{Row=3; Col=4; Value="1.876.22"; ForceText=$True},
    @{Row=3; Col=5; Value="  -0.61%  "; ForceText=$True},
    @{Row=4; Col=5; Value="  +0.08%  "; ForceText=$True},
    @{Row=5; Col=4; Value="0.7389"; ForceText=$True},
    @{Row=5; Col=5; Value="  -3.94%  "; ForceText=$True},
    @{Row=6; Col=4; Value="242.73"; ForceText=$True},
    @{Row=6; Col=5; Value="  +0.04%  "; ForceText=$True},
    @{Row=7; Col=4; Value="1.001"; ForceText=$True},
    @{Row=7; Col=5; Value="  +0.11%  "; ForceText=$True},
    @{Row=8; Col=5; Value="  +1.10%  "; ForceText=$True},
    @{Row=9; Col=4; Value="0.07199"; ForceText=$True},
    @{Row=9; Col=5; Value="  +0.47%  "; ForceText=$True},
    @{Row=10; Col=4; Value="24.60"; ForceText=$True},
    @{Row=10; Col=5; Value="  -4.04%  "; ForceText=$True},
    @{Row=11; Col=4; Value="0.08333"; ForceText=$True},
    @{Row=11; Col=5; Value="  -2.86%  "; ForceText=$True},
    @{Row=12; Col=4; Value="0.7496"; ForceText=$True},
    @{Row=13; Col=2; Value="WrappedEther"; ForceText=$False},
    @{Row=13; Col=3; Value="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; ForceText=$False},
    @{Row=13; Col=4; Value="1.894.24"; ForceText=$True},
    @{Row=13; Col=5; Value="  -1.21%  "; ForceText=$True},
    @{Row=14; Col=2; Value="Polkadot"; ForceText=$False},
    @{Row=14; Col=3; Value="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; ForceText=$False},
    @{Row=14; Col=4; Value="5.390"; ForceText=$True},
    @{Row=14; Col=5; Value="  +0.52%  "; ForceText=$True},
    @{Row=15; Col=4; Value="92.33"; ForceText=$True},
    @{Row=15; Col=5; Value="  -1.33%  "; ForceText=$True},
    @{Row=16; Col=4; Value="29.918.01"; ForceText=$True},
    @{Row=16; Col=5; Value="  +0.12%  "; ForceText=$True},
    @{Row=17; Col=4; Value="6.089"; ForceText=$True},
    @{Row=17; Col=5; Value="  -0.91%  "; ForceText=$True},
    @{Row=18; Col=4; Value="248.06"; ForceText=$True},
    @{Row=18; Col=5; Value="  +1.48%  "; ForceText=$True},
    @{Row=19; Col=4; Value="13.55"; ForceText=$True},
    @{Row=19; Col=5; Value="  -1.52%  "; ForceText=$True},
    @{Row=20; Col=4; Value="0.000007835"; ForceText=$True},
    @{Row=20; Col=5; Value="  +0.42%  "; ForceText=$True},
    @{Row=21; Col=2; Value="WrappedliquidstakedEther2.0"; ForceText=$False},
    @{Row=21; Col=3; Value="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; ForceText=$False},
    @{Row=21; Col=4; Value="2.141.90"; ForceText=$True},
    @{Row=21; Col=5; Value="  -1.32%  "; ForceText=$True},
    @{Row=22; Col=2; Value="Dai"; ForceText=$False},
    @{Row=22; Col=3; Value="https://coinranking.com/coin/MoTuySvg7+dai-dai"; ForceText=$False},
    @{Row=22; Col=4; Value="0.9992"; ForceText=$True},
    @{Row=22; Col=5; Value="  -0.05%  "; ForceText=$True},
    @{Row=23; Col=4; Value="8.005"; ForceText=$True},
    @{Row=23; Col=5; Value="  -0.05%  "; ForceText=$True},
    @{Row=24; Col=5; Value="  +0.00%  "; ForceText=$True},
    @{Row=25; Col=4; Value="0.1549"; ForceText=$True},
    @{Row=25; Col=5; Value="  -5.35%  "; ForceText=$True},
    @{Row=26; Col=4; Value="9.268"; ForceText=$True},
    @{Row=26; Col=5; Value="  -1.18%  "; ForceText=$True},
    @{Row=27; Col=4; Value="164.67"; ForceText=$True},
    @{Row=27; Col=5; Value="  +1.20%  "; ForceText=$True},
    @{Row=28; Col=4; Value="18.66"; ForceText=$True},
    @{Row=28; Col=5; Value="  -0.36%  "; ForceText=$True},
    @{Row=29; Col=4; Value="2.029"; ForceText=$True},
    @{Row=29; Col=5; Value="  -0.12%  "; ForceText=$True},
    @{Row=30; Col=4; Value="1.509"; ForceText=$True},
    @{Row=30; Col=5; Value="  +3.24%  "; ForceText=$True},
    @{Row=31; Col=4; Value="4.595"; ForceText=$True},
    @{Row=31; Col=5; Value="  +1.73%  "; ForceText=$True},
    @{Row=32; Col=4; Value="1.533"; ForceText=$True},
    @{Row=32; Col=5; Value="  -0.34%  "; ForceText=$True},
    @{Row=33; Col=4; Value="4.262"; ForceText=$True},
    @{Row=33; Col=5; Value="  +4.10%  "; ForceText=$True},
    @{Row=34; Col=4; Value="0.05318"; ForceText=$True},
    @{Row=34; Col=5; Value="  -2.55%  "; ForceText=$True},
    @{Row=35; Col=4; Value="1.234"; ForceText=$True},
    @{Row=35; Col=5; Value="  -0.46%  "; ForceText=$True},
    @{Row=36; Col=4; Value="0.7479"; ForceText=$True},
    @{Row=36; Col=5; Value="  +0.67%  "; ForceText=$True},
    @{Row=37; Col=4; Value="1.000"; ForceText=$True},
    @{Row=37; Col=5; Value="  -0.27%  "; ForceText=$True},
    @{Row=38; Col=5; Value="  -0.06%  "; ForceText=$True},
    @{Row=39; Col=4; Value="0.01964"; ForceText=$True},
    @{Row=39; Col=5; Value="  +0.56%  "; ForceText=$True},
    @{Row=40; Col=4; Value="2.755"; ForceText=$True},
    @{Row=40; Col=5; Value="  -1.02%  "; ForceText=$True},
    @{Row=41; Col=4; Value="0.4540"; ForceText=$True},
    @{Row=41; Col=5; Value="  +1.70%  "; ForceText=$True},
    @{Row=42; Col=2; Value="Maker"; ForceText=$False},
    @{Row=42; Col=3; Value="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; ForceText=$False},
    @{Row=42; Col=4; Value="1.106.62"; ForceText=$True},
    @{Row=42; Col=5; Value="  -0.25%  "; ForceText=$True},
    @{Row=43; Col=2; Value="FraxShare"; ForceText=$False},
    @{Row=43; Col=3; Value="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; ForceText=$False},
    @{Row=43; Col=4; Value="6.123"; ForceText=$True},
    @{Row=43; Col=5; Value="  +0.78%  "; ForceText=$True},
    @{Row=44; Col=4; Value="72.25"; ForceText=$True},
    @{Row=44; Col=5; Value="  -1.11%  "; ForceText=$True},
    @{Row=45; Col=4; Value="0.8594"; ForceText=$True},
    @{Row=45; Col=5; Value="  +0.86%  "; ForceText=$True},
    @{Row=46; Col=2; Value="PaxDollar"; ForceText=$False},
    @{Row=46; Col=3; Value="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; ForceText=$False},
    @{Row=46; Col=4; Value="1.003"; ForceText=$True},
    @{Row=46; Col=5; Value="  +0.24%  "; ForceText=$True},
    @{Row=47; Col=2; Value="Quant"; ForceText=$False},
    @{Row=47; Col=3; Value="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; ForceText=$False},
    @{Row=47; Col=4; Value="104.15"; ForceText=$True},
    @{Row=47; Col=5; Value="  +1.70%  "; ForceText=$True},
    @{Row=48; Col=4; Value="1.853"; ForceText=$True},
    @{Row=48; Col=5; Value="  -0.34%  "; ForceText=$True},
    @{Row=49; Col=4; Value="7.600"; ForceText=$True},
    @{Row=49; Col=5; Value="  -0.50%  "; ForceText=$True},
    @{Row=50; Col=4; Value="9.511"; ForceText=$True},
    @{Row=50; Col=5; Value="  -2.50%  "; ForceText=$True},
    @{Row=51; Col=4; Value="2.037.61"; ForceText=$True},
    @{Row=51; Col=5; Value="  -0.81%  "; ForceText=$True}
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
